$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# (e.g. "579.00" -> 579, "188.38" -> 188.38) and the formatting/leading
# zeros / trailing zeros would be lost.

$ws.Range("D2").Value = "66.841.20"
$ws.Range("E2").Value = "  +2.74%  "
$ws.Range("D3").Value = "3.436.53"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.00"
$ws.Range("E5").Value = "  +4.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.38"
$ws.Range("E6").Value = "  +8.29%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "3.428.24"
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.97"
$ws.Range("E12").Value = "  +6.41%  "
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.44"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "3.984.74"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.81"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").Value = "3.438.83"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "66.977.68"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.04"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "476.28"
$ws.Range("E22").Value = "  +3.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.40"
$ws.Range("E23").Value = "  +10.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.01"
$ws.Range("E24").Value = "  +20.03%  "
$ws.Range("E25").Value = "  +6.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.97"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.97"
$ws.Range("E27").Value = "  +3.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.92"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.97"
$ws.Range("E29").Value = "  +3.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.03"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.36"
$ws.Range("E31").Value = "  +12.86%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "599.26"
$ws.Range("E32").Value = "  +3.92%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.61"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.73"
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.111"
$ws.Range("E35").Value = "  +3.66%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.96"
$ws.Range("E38").Value = "  +3.85%  "
$ws.Range("E39").Value = "  +4.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.46"
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("D41").Value = "0.0₃0749"
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").Value = "3.195.36"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.89"
$ws.Range("E43").Value = "  +5.34%  "
$ws.Range("E44").Value = "  +3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.57"
$ws.Range("E45").Value = "  +5.15%  "
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("E47").Value = "  +20.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.134"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.19"
$ws.Range("E50").Value = "  +5.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.56"
$ws.Range("E51").Value = "  +3.33%  "
